# Auto-generated edit script applying the Exodus_Profits value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1015.7451
$ws.Range("I15").Value = 1015.7451
$ws.Range("K15").Value = 3047.2353
$ws.Range("M15").Value = -2878.2353
$ws.Range("H70").Value = 2545.2
$ws.Range("I70").Value = 1010
$ws.Range("J70").Value = 2929
$ws.Range("K70").Value = 3030
$ws.Range("L70").Value = 8787
$ws.Range("M70").Value = -2760
$ws.Range("N70").Value = -9327
$ws.Range("H73").Value = 2545.2
$ws.Range("I73").Value = 1010
$ws.Range("J73").Value = 2929
$ws.Range("K73").Value = 3030
$ws.Range("L73").Value = 8787
$ws.Range("M73").Value = -2094
$ws.Range("N73").Value = -10659
$ws.Range("H74").Value = 5414.5713
$ws.Range("I74").Value = 5289.9
$ws.Range("K74").Value = 5289.9
$ws.Range("M74").Value = -4353.9
$ws.Range("H77").Value = 5414.5713
$ws.Range("I77").Value = 5289.9
$ws.Range("K77").Value = 26449.5
$ws.Range("M77").Value = -21769.5
$ws.Range("H98").Value = 1451.125
$ws.Range("I98").Value = 1119.7142
$ws.Range("J98").Value = 1708.8889
$ws.Range("K98").Value = 1119.7142
$ws.Range("L98").Value = 1708.8889
$ws.Range("M98").Value = 378.2858000000001
$ws.Range("N98").Value = -4704.8889
$ws.Range("H107").Value = 132.58824
$ws.Range("I107").Value = 143.21428
$ws.Range("J107").Value = 83
$ws.Range("K107").Value = 143.21428
$ws.Range("L107").Value = 83
$ws.Range("M107").Value = 1776.78572
$ws.Range("N107").Value = -3923
$ws.Range("H122").Value = 1451.125
$ws.Range("I122").Value = 1119.7142
$ws.Range("J122").Value = 1708.8889
$ws.Range("K122").Value = 3359.1426
$ws.Range("L122").Value = 5126.6667
$ws.Range("M122").Value = -909.1425999999997
$ws.Range("N122").Value = -10026.6667
$ws.Range("H125").Value = 4535.273
$ws.Range("I125").Value = 4158.8
$ws.Range("J125").Value = 4849
$ws.Range("K125").Value = 37429.2
$ws.Range("L125").Value = 43641
$ws.Range("M125").Value = -34969.2
$ws.Range("N125").Value = -48561
$ws.Range("H133").Value = 76650.60000000001
$ws.Range("J133").Value = 76650.60000000001
$ws.Range("L133").Value = 76650.60000000001
$ws.Range("N133").Value = -86770.60000000001
$ws.Range("H134").Value = 94216.664
$ws.Range("J134").Value = 94216.664
$ws.Range("L134").Value = 94216.664
$ws.Range("N134").Value = -104356.664
$ws.Range("H136").Value = 96495.836
$ws.Range("J136").Value = 96495.836
$ws.Range("L136").Value = 96495.836
$ws.Range("N136").Value = -106695.836
$ws.Range("H137").Value = 286204.6
$ws.Range("I137").Value = 1702.5667
$ws.Range("J137").Value = 692636.1
$ws.Range("K137").Value = 5107.7001
$ws.Range("L137").Value = 2077908.3
$ws.Range("M137").Value = -2557.7001
$ws.Range("N137").Value = -2083008.3
$ws.Range("H139").Value = 99406
$ws.Range("J139").Value = 99406
$ws.Range("L139").Value = 99406
$ws.Range("N139").Value = -109686
$ws.Range("H140").Value = 89423.57000000001
$ws.Range("J140").Value = 89423.57000000001
$ws.Range("L140").Value = 89423.57000000001
$ws.Range("N140").Value = -99783.57000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2296.2856
$ws.Range("I132").Value = 1945.0952
$ws.Range("J132").Value = 3349.8572
$ws.Range("K132").Value = 5835.2856
$ws.Range("L132").Value = 10049.5716
$ws.Range("M132").Value = -3305.2856
$ws.Range("N132").Value = -15109.5716

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 3284.7778
$ws.Range("I54").Value = 1820.375
$ws.Range("K54").Value = 1820.375
$ws.Range("M54").Value = -1336.375
$ws.Range("H107").Value = 2223.5454
$ws.Range("I107").Value = 2001
$ws.Range("J107").Value = 2817
$ws.Range("K107").Value = 2001
$ws.Range("L107").Value = 2817
$ws.Range("M107").Value = -81
$ws.Range("N107").Value = -6657
$ws.Range("H132").Value = 46554.71
$ws.Range("J132").Value = 46554.71
$ws.Range("L132").Value = 46554.71
$ws.Range("N132").Value = -56674.71
$ws.Range("H135").Value = 104024.336
$ws.Range("J135").Value = 104024.336
$ws.Range("L135").Value = 104024.336
$ws.Range("N135").Value = -114164.336
$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279
$ws.Range("H140").Value = 43476
$ws.Range("J140").Value = 43499.406
$ws.Range("L140").Value = 43499.406
$ws.Range("N140").Value = -53859.406

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 18958.5
$ws.Range("I94").Value = 22627.8
$ws.Range("J94").Value = 612
$ws.Range("K94").Value = 22627.8
$ws.Range("L94").Value = 612
$ws.Range("M94").Value = -22176.8
$ws.Range("N94").Value = -1514
$ws.Range("H138").Value = 102243.43
$ws.Range("J138").Value = 108332.5
$ws.Range("L138").Value = 108332.5
$ws.Range("N138").Value = -118612.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 286
$ws.Range("I92").Value = 238.6
$ws.Range("K92").Value = 715.8
$ws.Range("M92").Value = 532.2
$ws.Range("M98").ClearContents()
$ws.Range("H98").Value = 9999.5
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 9999.5
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 29998.5
$ws.Range("N98").Value = -32994.5
$ws.Range("N107").ClearContents()
$ws.Range("H107").Value = 454.42856
$ws.Range("I107").Value = 454.42856
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1363.28568
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 556.71432
$ws.Range("H113").Value = 36230.395
$ws.Range("I113").Value = 322
$ws.Range("J113").Value = 56179.5
$ws.Range("K113").Value = 966
$ws.Range("L113").Value = 168538.5
$ws.Range("M113").Value = 1204
$ws.Range("N113").Value = -172878.5
$ws.Range("N122").ClearContents()
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6550

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 8999
$ws.Range("J4").Value = 8999
$ws.Range("L4").Value = 8999
$ws.Range("N4").Value = -9223
$ws.Range("H70").Value = 353819.25
$ws.Range("I70").Value = 171759
$ws.Range("J70").Value = 900000
$ws.Range("K70").Value = 171759
$ws.Range("L70").Value = 900000
$ws.Range("M70").Value = -171489
$ws.Range("N70").Value = -900540
$ws.Range("H73").Value = 353819.25
$ws.Range("I73").Value = 171759
$ws.Range("J73").Value = 900000
$ws.Range("K73").Value = 171759
$ws.Range("L73").Value = 900000
$ws.Range("M73").Value = -170823
$ws.Range("N73").Value = -901872
$ws.Range("H93").Value = 20596.9
$ws.Range("J93").Value = 20596.9
$ws.Range("L93").Value = 20596.9
$ws.Range("N93").Value = -24340.9
$ws.Range("H122").Value = 360406.22
$ws.Range("I122").Value = 628084.5600000001
$ws.Range("J122").Value = 3501.75
$ws.Range("K122").Value = 1884253.68
$ws.Range("L122").Value = 10505.25
$ws.Range("M122").Value = -1881803.68
$ws.Range("N122").Value = -15405.25
$ws.Range("H132").Value = 4175.1177
$ws.Range("I132").Value = 2760.875
$ws.Range("J132").Value = 5432.222
$ws.Range("K132").Value = 8282.625
$ws.Range("L132").Value = 16296.666
$ws.Range("M132").Value = -5752.625
$ws.Range("N132").Value = -21356.666
$ws.Range("H135").Value = 53519.523
$ws.Range("J135").Value = 53519.523
$ws.Range("L135").Value = 53519.523
$ws.Range("N135").Value = -63659.523
$ws.Range("H140").Value = 97214
$ws.Range("J140").Value = 98053.25
$ws.Range("L140").Value = 98053.25
$ws.Range("N140").Value = -108413.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4100
$ws.Range("I61").Value = 4200
$ws.Range("K61").Value = 4200
$ws.Range("M61").Value = -3998
$ws.Range("H113").Value = 4100
$ws.Range("I113").Value = 4200
$ws.Range("K113").Value = 4200
$ws.Range("M113").Value = -2030
$ws.Range("H121").Value = 44499.332
$ws.Range("J121").Value = 44499.332
$ws.Range("L121").Value = 44499.332
$ws.Range("N121").Value = -47993.332
$ws.Range("H122").Value = 60004540
$ws.Range("J122").Value = 28576492
$ws.Range("L122").Value = 85729476
$ws.Range("N122").Value = -85734376

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3560.5557
$ws.Range("I126").Value = 2591
$ws.Range("K126").Value = 7773
$ws.Range("M126").Value = -5303
$ws.Range("H136").Value = 1968
$ws.Range("I136").Value = 1767.1
$ws.Range("K136").Value = 5301.299999999999
$ws.Range("M136").Value = -2751.299999999999
